$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 149); this shifts nothing below it
# and reduces the sheet dimension from AC149 to AC148.
$ws.Rows.Item(149).Delete()

# Row 146: new match data (replacing the old row 146 content)
$ws.Cells.Item(146, 2).Value = 6769308          # B - id
$ws.Cells.Item(146, 5).Value = 45394.5          # E - Date
$ws.Cells.Item(146, 6).Value = "NK Varazdin"    # F - HomeTeam
$ws.Cells.Item(146, 7).Value = "Slaven Belupo"  # G - AwayTeam
$ws.Cells.Item(146, 8).Value = 1                # H - FTHG
$ws.Cells.Item(146, 9).Value = 3                # I - FTAG
$ws.Cells.Item(146, 10).Value = "A"             # J - FTR
$ws.Cells.Item(146, 11).Value = 4               # K - oddH_op
$ws.Cells.Item(146, 12).Value = 3.4             # L - oddD_op
$ws.Cells.Item(146, 13).Value = 1.95            # M - oddA_op
$ws.Cells.Item(146, 14).Value = 4.75            # N - oddH
$ws.Cells.Item(146, 15).Value = 3.5             # O - oddD
$ws.Cells.Item(146, 16).Value = 1.8             # P - oddA
$ws.Cells.Item(146, 17).Value = 0.5             # Q - Ah
$ws.Cells.Item(146, 18).Value = 2.05            # R - oddAHH
$ws.Cells.Item(146, 19).Value = 1.8             # S - oddAHA
$ws.Cells.Item(146, 20).Value = 2.5             # T - AhOU
$ws.Cells.Item(146, 21).Value = 1.975           # U - oddAHOver
$ws.Cells.Item(146, 22).Value = 1.875           # V - oddAHUnder
$ws.Cells.Item(146, 23).Value = -1              # W - PLH
$ws.Cells.Item(146, 24).Value = -1              # X - PLD
$ws.Cells.Item(146, 25).Value = 0.8             # Y - PLA
$ws.Cells.Item(146, 26).Value = -1              # Z - PL_Ahh
$ws.Cells.Item(146, 27).Value = 0.8             # AA - PL_Aha
$ws.Cells.Item(146, 28).Value = 0.9750000000000001 # AB - PL_AhOver
$ws.Cells.Item(146, 29).Value = -1              # AC - PL_AhUnder

# Row 147: now holds what used to be row 148's match, with a few odds refreshed
$ws.Cells.Item(147, 2).Value = 6923266          # B - id
$ws.Cells.Item(147, 5).Value = 45396.47916666666 # E - Date
$ws.Cells.Item(147, 6).Value = "Hajduk Split"   # F - HomeTeam
$ws.Cells.Item(147, 7).Value = "NK Osijek"      # G - AwayTeam
$ws.Cells.Item(147, 11).Value = 1.615           # K - oddH_op
$ws.Cells.Item(147, 12).Value = 3.5             # L - oddD_op
$ws.Cells.Item(147, 13).Value = 6               # M - oddA_op
$ws.Cells.Item(147, 14).Value = 1.6             # N - oddH
$ws.Cells.Item(147, 15).Value = 3.6             # O - oddD
$ws.Cells.Item(147, 16).Value = 6               # P - oddA
$ws.Cells.Item(147, 17).Value = -1              # Q - Ah
$ws.Cells.Item(147, 18).Value = 2.1             # R - oddAHH
$ws.Cells.Item(147, 19).Value = 1.775           # S - oddAHA
$ws.Cells.Item(147, 20).Value = 2.5             # T - AhOU
$ws.Cells.Item(147, 21).Value = 2.05            # U - oddAHOver
$ws.Cells.Item(147, 22).Value = 1.8             # V - oddAHUnder

# Row 148: now holds what used to be row 149's match, with a few odds refreshed
$ws.Cells.Item(148, 2).Value = 6788944          # B - id
$ws.Cells.Item(148, 5).Value = 45396.58333333334 # E - Date
$ws.Cells.Item(148, 6).Value = "Istra 1961"     # F - HomeTeam
$ws.Cells.Item(148, 7).Value = "HNK Rijeka"     # G - AwayTeam
$ws.Cells.Item(148, 11).Value = 5.5             # K - oddH_op
$ws.Cells.Item(148, 12).Value = 3.6             # L - oddD_op
$ws.Cells.Item(148, 13).Value = 1.615           # M - oddA_op
$ws.Cells.Item(148, 14).Value = 5.75            # N - oddH
$ws.Cells.Item(148, 15).Value = 3.75            # O - oddD
$ws.Cells.Item(148, 16).Value = 1.571           # P - oddA
$ws.Cells.Item(148, 17).Value = 1               # Q - Ah
$ws.Cells.Item(148, 18).Value = 1.775           # R - oddAHH
$ws.Cells.Item(148, 19).Value = 2.1             # S - oddAHA
$ws.Cells.Item(148, 21).Value = 1.875           # U - oddAHOver
$ws.Cells.Item(148, 22).Value = 1.975           # V - oddAHUnder
